# Update Betfair back/lay odds values to match the 2026-02-06 refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Cells.Item(2, 8).Value = 2.1
$ws.Cells.Item(2, 19).Value = 2.76
$ws.Cells.Item(2, 23).Value = 1.39
$ws.Cells.Item(2, 27).Value = 27

$ws.Cells.Item(3, 19).Value = 2.88

$ws.Cells.Item(4, 6).Value = 4.9
$ws.Cells.Item(4, 7).Value = 7.6
$ws.Cells.Item(4, 10).Value = 3.55
$ws.Cells.Item(4, 11).Value = 4.2
$ws.Cells.Item(4, 16).Value = 1.64
$ws.Cells.Item(4, 17).Value = 1.89

$ws.Cells.Item(5, 6).Value = 1.63
$ws.Cells.Item(5, 7).Value = 1.74
$ws.Cells.Item(5, 8).Value = 5.2
$ws.Cells.Item(5, 9).Value = 6.2
$ws.Cells.Item(5, 11).Value = 4.7
$ws.Cells.Item(5, 12).Value = 1.3
$ws.Cells.Item(5, 13).Value = 1.06
$ws.Cells.Item(5, 14).Value = 4
$ws.Cells.Item(5, 15).Value = 1.27
$ws.Cells.Item(5, 16).Value = 2.06
$ws.Cells.Item(5, 17).Value = 1.79
$ws.Cells.Item(5, 18).Value = 1.41
$ws.Cells.Item(5, 19).Value = 3.05
$ws.Cells.Item(5, 20).Value = 1.83
$ws.Cells.Item(5, 21).Value = 2.02
$ws.Cells.Item(5, 22).Value = 1.19
$ws.Cells.Item(5, 23).Value = 2.34
$ws.Cells.Item(5, 24).Value = 22
$ws.Cells.Item(5, 25).Value = 19.5
$ws.Cells.Item(5, 26).Value = 46
$ws.Cells.Item(5, 27).Value = 170
$ws.Cells.Item(5, 28).Value = 9.199999999999999
$ws.Cells.Item(5, 29).Value = 10.5
$ws.Cells.Item(5, 30).Value = 21
$ws.Cells.Item(5, 31).Value = 75
$ws.Cells.Item(5, 32).Value = 10.5
$ws.Cells.Item(5, 33).Value = 10.5
$ws.Cells.Item(5, 34).Value = 23
$ws.Cells.Item(5, 35).Value = 75
$ws.Cells.Item(5, 36).Value = 16.5
$ws.Cells.Item(5, 37).Value = 17.5
$ws.Cells.Item(5, 38).Value = 36
$ws.Cells.Item(5, 39).Value = 140
$ws.Cells.Item(5, 40).Value = 9.199999999999999
$ws.Cells.Item(5, 41).Value = 1000

$ws.Cells.Item(6, 6).Value = 3.15
$ws.Cells.Item(6, 7).Value = 3.65
$ws.Cells.Item(6, 8).Value = 2.68
$ws.Cells.Item(6, 9).Value = 3.1
$ws.Cells.Item(6, 12).Value = 1.57
$ws.Cells.Item(6, 13).Value = 1.14
$ws.Cells.Item(6, 14).Value = 2.2
$ws.Cells.Item(6, 15).Value = 1.67
$ws.Cells.Item(6, 17).Value = 3.05
$ws.Cells.Item(6, 18).Value = 1.13
$ws.Cells.Item(6, 19).Value = 6.8
$ws.Cells.Item(6, 20).Value = 2.28
$ws.Cells.Item(6, 21).Value = 1.63
$ws.Cells.Item(6, 22).Value = 1.49
$ws.Cells.Item(6, 23).Value = 1.37
$ws.Cells.Item(6, 24).Value = 8.199999999999999
$ws.Cells.Item(6, 25).Value = 8.800000000000001
$ws.Cells.Item(6, 26).Value = 20
$ws.Cells.Item(6, 27).Value = 65
$ws.Cells.Item(6, 28).Value = 10
$ws.Cells.Item(6, 29).Value = 8.4
$ws.Cells.Item(6, 30).Value = 17.5
$ws.Cells.Item(6, 31).Value = 60
$ws.Cells.Item(6, 32).Value = 26
$ws.Cells.Item(6, 33).Value = 20
$ws.Cells.Item(6, 34).Value = 970
$ws.Cells.Item(6, 35).Value = 110
$ws.Cells.Item(6, 36).Value = 90
$ws.Cells.Item(6, 37).Value = 80
$ws.Cells.Item(6, 38).Value = 130
$ws.Cells.Item(6, 39).Value = 330
$ws.Cells.Item(6, 40).Value = 120
$ws.Cells.Item(6, 41).Value = 80

$ws.Cells.Item(7, 6).Value = 2.22
$ws.Cells.Item(7, 10).Value = 2.92
$ws.Cells.Item(7, 12).Value = 1.51
$ws.Cells.Item(7, 13).Value = 1.13
$ws.Cells.Item(7, 14).Value = 2.38
$ws.Cells.Item(7, 15).Value = 1.59
$ws.Cells.Item(7, 17).Value = 2.82
$ws.Cells.Item(7, 18).Value = 1.18
$ws.Cells.Item(7, 19).Value = 5.7
$ws.Cells.Item(7, 20).Value = 2.2
$ws.Cells.Item(7, 21).Value = 1.68
$ws.Cells.Item(7, 22).Value = 1.28
$ws.Cells.Item(7, 23).Value = 1.69
$ws.Cells.Item(7, 24).Value = 9.4
$ws.Cells.Item(7, 25).Value = 11
$ws.Cells.Item(7, 26).Value = 29
$ws.Cells.Item(7, 27).Value = 130
$ws.Cells.Item(7, 28).Value = 7
$ws.Cells.Item(7, 29).Value = 7.4
$ws.Cells.Item(7, 30).Value = 19.5
$ws.Cells.Item(7, 31).Value = 80
$ws.Cells.Item(7, 32).Value = 13.5
$ws.Cells.Item(7, 33).Value = 13
$ws.Cells.Item(7, 34).Value = 32
$ws.Cells.Item(7, 35).Value = 130
$ws.Cells.Item(7, 36).Value = 36
$ws.Cells.Item(7, 37).Value = 38
$ws.Cells.Item(7, 38).Value = 75
$ws.Cells.Item(7, 39).Value = 280
$ws.Cells.Item(7, 40).Value = 46
$ws.Cells.Item(7, 41).Value = 150

$ws.Cells.Item(8, 6).Value = 1.88
$ws.Cells.Item(8, 7).Value = 2.02
$ws.Cells.Item(8, 9).Value = 4.9
$ws.Cells.Item(8, 10).Value = 3.5
$ws.Cells.Item(8, 12).Value = 1.37
$ws.Cells.Item(8, 13).Value = 1.07
$ws.Cells.Item(8, 14).Value = 3.4
$ws.Cells.Item(8, 15).Value = 1.35
$ws.Cells.Item(8, 16).Value = 1.82
$ws.Cells.Item(8, 17).Value = 1.93
$ws.Cells.Item(8, 18).Value = 1.31
$ws.Cells.Item(8, 19).Value = 3.65
$ws.Cells.Item(8, 20).Value = 1.86
$ws.Cells.Item(8, 21).Value = 1.98
$ws.Cells.Item(8, 22).Value = 1.25
$ws.Cells.Item(8, 23).Value = 1.98
$ws.Cells.Item(8, 24).Value = 15.5
$ws.Cells.Item(8, 25).Value = 17.5
$ws.Cells.Item(8, 26).Value = 40
$ws.Cells.Item(8, 27).Value = 130
$ws.Cells.Item(8, 28).Value = 9.800000000000001
$ws.Cells.Item(8, 29).Value = 9.4
$ws.Cells.Item(8, 30).Value = 22
$ws.Cells.Item(8, 31).Value = 75
$ws.Cells.Item(8, 32).Value = 13.5
$ws.Cells.Item(8, 33).Value = 12
$ws.Cells.Item(8, 34).Value = 23
$ws.Cells.Item(8, 35).Value = 85
$ws.Cells.Item(8, 36).Value = 27
$ws.Cells.Item(8, 37).Value = 25
$ws.Cells.Item(8, 38).Value = 46
$ws.Cells.Item(8, 39).Value = 140
$ws.Cells.Item(8, 40).Value = 18
$ws.Cells.Item(8, 41).Value = 90

$ws.Cells.Item(9, 6).Value = 1.52
$ws.Cells.Item(9, 9).Value = 13
$ws.Cells.Item(9, 10).Value = 3.9
$ws.Cells.Item(9, 11).Value = 4.5
$ws.Cells.Item(9, 16).Value = 1.54
$ws.Cells.Item(9, 17).Value = 2.08

$ws.Cells.Item(10, 6).Value = 1.64
$ws.Cells.Item(10, 7).Value = 1.79
$ws.Cells.Item(10, 8).Value = 4.5
$ws.Cells.Item(10, 11).Value = 5.1

$ws.Cells.Item(12, 9).Value = 1000

$ws.Cells.Item(13, 7).Value = 1.81
$ws.Cells.Item(13, 9).Value = 7
$ws.Cells.Item(13, 16).Value = 1.68
$ws.Cells.Item(13, 17).Value = 2.18

$ws.Cells.Item(14, 6).Value = 2.16
$ws.Cells.Item(14, 10).Value = 3.65
$ws.Cells.Item(14, 11).Value = 3.7
$ws.Cells.Item(14, 16).Value = 1.98
$ws.Cells.Item(14, 18).Value = 1.38
$ws.Cells.Item(14, 26).Value = 27
$ws.Cells.Item(14, 29).Value = 8.199999999999999

$ws.Cells.Item(15, 6).Value = 4.7
$ws.Cells.Item(15, 11).Value = 4.2
$ws.Cells.Item(15, 20).Value = 1.8
$ws.Cells.Item(15, 24).Value = 17.5

$ws.Cells.Item(16, 10).Value = 3.15
$ws.Cells.Item(16, 11).Value = 3.2
$ws.Cells.Item(16, 24).Value = 9.4
$ws.Cells.Item(16, 29).Value = 7

$ws.Cells.Item(17, 18).Value = 1.3
$ws.Cells.Item(17, 20).Value = 1.91
$ws.Cells.Item(17, 41).Value = 60

$ws.Cells.Item(18, 7).Value = 2.08
$ws.Cells.Item(18, 31).Value = 70
